$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("machine")

# Update the calendar id in H2 from "CAL-PADRAO-5x8" to "CAL-24x5"
$ws.Range("H2").Value = "CAL-24x5"

# Update the active selection from H1 to I1
$ws.Range("I1").Select()
